$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: B2 (zh-cn status), C2 (de-de status), D2 (Latest Handoff Date)
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"
$overview.Range("D2").Value = "2016-28-20 06:28:26"

# zh-cn sheet: C2 (Status), E2 (Latest Handoff Datetime)
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("E2").Value = "2016-03-20 06:28:23"

# de-de sheet: C2 (Status), E2 (Latest Handoff Datetime)
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("E2").Value = "2016-03-20 06:28:26"
